# Automatic update 2025-07-08 11:40:09
# Inserts a new advisor row ("ALCIVAR BUSTAMANTE ERNESTO EDUARDO" under
# "OFICINA-CATAECSA") at row 241 in both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, pushing the existing rows 241-281 down to
# 242-282, and refreshes the trailing summary row text/count that
# referenced the old total of 279 records.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" (columns A:R, summary row shows "N de 279")
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Push rows 241..281 down to 242..282, leaving a blank row 241 behind.
$ws1.Rows.Item(241).Insert()

# Populate the newly inserted row with the new advisor entry.
$ws1.Cells.Item(241, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(241, 2).Value = "ALCIVAR BUSTAMANTE ERNESTO EDUARDO"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(241, $c).Value = 0
}

# The trailing "N de 279" counters now cover one more record ("N de 280").
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(282, $c)
    $cell.Value = ($cell.Value2 -replace "de 279", "de 280")
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" (columns A:G, summary row is numeric totals)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(241).Insert()

$ws2.Cells.Item(241, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(241, 2).Value = "ALCIVAR BUSTAMANTE ERNESTO EDUARDO"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(241, $c).Value = 0
}
